# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header "Subj" counts): update submax values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON"): B2/C2/D2 values were dropped (now blank), E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -4.2238028298040078

# Row 3 ("STR"): B3 dropped, C3 updated, D3 newly added, E3 updated
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 0.5374288134181171
$ws.Range("D3").Value = -7.3463355315719454
$ws.Range("E3").Value = 9.8447801482831068

# Selection now highlights only the edited block B1:E3
$ws.Range("B1:E3").Select()
